$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.083.23'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.926.80'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.59%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.07'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.10'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +0.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.922.84'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.82'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.68%  '
$ws.Range('E11').Value = '  -0.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.443'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.49%  '
$ws.Range('E13').Value = '  +1.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.70'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.45%  '
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.410.22'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.005.16'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.74'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.925.90'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '431.74'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.51'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.37%  '
$ws.Range('E22').Value = '  +1.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.10'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.28%  '
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.84'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.74%  '
$ws.Range('E26').Value = '  +1.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.14'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.90%  '
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('E29').Value = '  +6.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.62'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.11'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.58'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.108'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0863'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.88%  '
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.63'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.08'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '49.98'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('E40').Value = '  -0.74%  '
$ws.Range('E41').Value = '  -1.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.63'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.289'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '39.77'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '383.12'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.75%  '
$ws.Range('E46').Value = '  +0.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.705.69'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.73'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.23'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.107'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.23%  '
